$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 10 down to row 11 (same column styles), then set the new values
$ws.Range("A10:C10").Copy()
$ws.Range("A11:C11").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A11").Value = 43903
$ws.Range("B11").Value = 1.5
$ws.Range("C11").Value = "Entrevue avec le chef de projet pour savoir où j'en suis et les problèmes rencontré"

$ws.Rows.Item(11).RowHeight = 30

$ws.Range("B11").Select()
